# Apply Bitcoin price tracker update: shift current prices into the
# "Old Prices" column (D) and record freshly scraped prices (B) and their
# Euro conversions (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Binance): old "current price" becomes the new "Old Prices" entry.
$ws.Range("D2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = 46821.07
$ws.Range("C2").Value2 = 41562.910000000003

# Row 3 (KuCoin)
$ws.Range("D3").Value2 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = 46874.65
$ws.Range("C3").Value2 = 41610.47

# Row 4 (Coinbase)
$ws.Range("D4").Value2 = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = 46827.42
$ws.Range("C4").Value2 = 41568.550000000003

# Row 5 (CMC)
$ws.Range("D5").Value2 = $ws.Range("B5").Value2
$ws.Range("B5").Value2 = 46534.879999999997
$ws.Range("C5").Value2 = 41308.86

# Row 6 (CoinGecko) holds its scraped "price" columns as literal text
# (trailing non-breaking space, straight off the web page), so force the
# Text number format before assigning or the numeric-looking string gets
# silently reinterpreted as a number.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "46616.25 "
$ws.Range("D6").Style = $ws.Range("A1").Style

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "46662.41 "
$ws.Range("B6").Style = $ws.Range("A1").Style

$ws.Range("C6").Value2 = 41422.07

$wb.Save()
